{"js": "// Add a new paragraph \"Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement }}\"\n// right after the \"N\u00b0 SIRET : {{ etablissement.siret }}\" paragraph, in the\n// TIAC \"etablissement\" block \u2014 matching the existing paragraph's style and\n// run formatting (Corpsdetexte, Calibri 10pt).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet siretParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  // Avoid relying on the exact punctuation/whitespace (the source uses a\n  // non-breaking space before the colon), \"SIRET\" alone is unambiguous here.\n  if (paragraphs.items[i].text.indexOf(\"SIRET\") !== -1) {\n    siretParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!siretParagraph) {\n  throw new Error('Could not locate the \"N\u00b0 SIRET\" paragraph.');\n}\n\n// Formatting shared by the paragraph mark and by the surrounding runs in\n// this block of the template (same as the \"N\u00b0 SIRET\" paragraph above it).\nconst fullRunFonts =\n  '<w:rFonts w:cs=\"Calibri\" w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>';\nconst midRunFonts =\n  '<w:rFonts w:cs=\"\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>';\nconst sizeTags = '<w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/>';\n\nconst newParagraphXml =\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"Corpsdetexte\"/><w:rPr>' +\n  '<w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>' +\n  sizeTags +\n  \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" +\n  fullRunFonts +\n  sizeTags +\n  \"</w:rPr><w:t>Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement</w:t></w:r>\" +\n  \"<w:r><w:rPr>\" +\n  midRunFonts +\n  sizeTags +\n  '</w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:rPr>\" +\n  fullRunFonts +\n  sizeTags +\n  \"</w:rPr><w:t>}}</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst insertionPoint = siretParagraph.getRange(Word.RangeLocation.end);\ninsertionPoint.insertOoxml(flatOpcPackage, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add \"Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement }}\" as a new\n# paragraph right after the \"N\u00b0 SIRET : {{ etablissement.siret }}\"\n# paragraph, matching that paragraph's style/formatting.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.Text = \"N\u00b0 SIRET :\"\n$rng.Find.Execute() | Out-Null\n\n$siretParagraph = $rng.Paragraphs(1)\n$insertionPoint = $siretParagraph.Range\n$insertionPoint.Collapse(0)\n$insertionPoint.InsertParagraphAfter()\n$insertionPoint.InsertAfter(\"Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement }}\")\n"}
